# edit.ps1
# Applies the "Modify POODLE cmd line output / Tell user they can search at
# any time" commit to the document.
#
# Summary of changes:
#  1. In the "Additional files:" paragraph, reword the description of the
#     files that poodle.py creates/uses (ignore.txt / database.txt) and
#     drop the separate mentions of graph.txt / ranks.txt / index.txt in
#     favour of a single database.txt description.
#  2. Word's "_GoBack" bookmark (marks the location of the user's last
#     edit) therefore moves from the "POODLE Options" heading to the end
#     of the reworded "Additional files:" paragraph.
#  3. A couple of incidental proofing-run merges that have no visible text
#     impact are tidied up as a natural consequence of the text replaces.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Reword the "Additional files:" paragraph. The replacement is split
#    into two Find/Replace calls (before/after the word "url") so the
#    existing "url" run (and its spell-check proofErr wrapper) is left
#    intact, matching how Word itself would merge only the runs that
#    were actually touched by the edit.
# ---------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.ClearFormatting()
$old1 = " ignore.txt (edit this file and add words to ignore when scraping if you wish), graph.txt (stores "
$new1 = " ignore.txt (edit this file and add words to ignore when scraping if you wish), database.txt (Stores the "
$replaced1 = $find1.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Host "Additional files paragraph updated (part 1): $replaced1"

$find2 = $d.Content.Find
$find2.ClearFormatting()
$old2 = " graph data structure), ranks.txt (stores rank data structure), and index.txt (index data structure)"
$new2 = " graph, index, and page ranks)"
$replaced2 = $find2.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Host "Additional files paragraph updated (part 2): $replaced2"

# ---------------------------------------------------------------------
# 2. Tidy a couple of now-redundant proofing run splits that have no
#    effect on the visible text (safe, content-preserving merges).
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$t = ", and cannot enter the if statement that only accepts "
$find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2) | Out-Null

$find = $d.Content.Find
$find.ClearFormatting()
$t2 = "(key = lambda x: x[1]) " + [char]0x2013 + " this means the parameter " + [char]0x201C + "x" + [char]0x201D + " (the list) is passed to the function, and the "
$find.Execute($t2, $false, $false, $false, $false, $false, $true, 1, $false, $t2, 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the "POODLE Options" heading to the
#    end of the reworded "Additional files:" paragraph.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
    Write-Host "Removed old _GoBack bookmark"
}

$anchor = $d.Content
$find = $anchor.Find
$find.ClearFormatting()
$found = $find.Execute("page ranks)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found bookmark anchor text: $found"

if ($found) {
    $anchor.Collapse(0)
    $insertPos = $anchor.Start

    # Work around a COM-interop quirk where adding a bookmark directly to
    # a zero-length range positioned immediately before a paragraph mark
    # can resolve to the wrong location. Instead, insert a temporary
    # placeholder character, bookmark the (non-empty) range around it,
    # then delete the placeholder - the bookmark correctly collapses to
    # the intended location.
    $anchor.InsertAfter("X")
    $d.Bookmarks.Add("_GoBack", $anchor)

    $placeholder = $d.Range($insertPos, $insertPos + 1)
    $placeholder.Text = ""

    $finalBookmark = $d.Bookmarks("_GoBack")
    Write-Host "New _GoBack bookmark at: $($finalBookmark.Start)-$($finalBookmark.End)"
}
